# Updates the crypto price table (Sheet1) cell-by-cell to match the latest
# coinranking.com snapshot referenced in the commit message. Two coin pairs
# (rows 14/15 and rows 43/44) swapped rank order, so Coin/Link (B/C) are
# rewritten for those rows too, not just Price/Volume (D/E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: cell reference, new value, and whether the value is a
# numeric-looking string that must be kept as TEXT (the sheet stores every
# Price/Volume cell as a string, e.g. "25.50", "586.01", "0.0000171" -
# some of those would otherwise be auto-converted to numbers by Excel and
# lose their exact formatting, e.g. "25.50" -> 25.5).
$updates = @(
    @{ Cell = 'D2'; Value = '67.221.86'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  +0.60%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '2.493.37'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  +0.47%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '586.01'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  +0.20%  '; ForceText = $false }
    @{ Cell = 'D6'; Value = '172.74'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  +2.64%  '; ForceText = $false }
    @{ Cell = 'E7'; Value = '  -0.10%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '0.514'; ForceText = $true }
    @{ Cell = 'E8'; Value = '  -0.41%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '2.488.92'; ForceText = $false }
    @{ Cell = 'E9'; Value = '  +0.26%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '0.137'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  +0.88%  '; ForceText = $false }
    @{ Cell = 'E11'; Value = '  -0.12%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '4.94'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  -0.18%  '; ForceText = $false }
    @{ Cell = 'D13'; Value = '0.333'; ForceText = $true }
    @{ Cell = 'E13'; Value = '  -1.02%  '; ForceText = $false }
    @{ Cell = 'B14'; Value = 'Avalanche'; ForceText = $false }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; ForceText = $false }
    @{ Cell = 'D14'; Value = '25.50'; ForceText = $true }
    @{ Cell = 'E14'; Value = '  -1.54%  '; ForceText = $false }
    @{ Cell = 'B15'; Value = 'WrappedliquidstakedEther2.0'; ForceText = $false }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; ForceText = $false }
    @{ Cell = 'D15'; Value = '2.919.62'; ForceText = $false }
    @{ Cell = 'E15'; Value = '  +1.21%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '66.943.23'; ForceText = $false }
    @{ Cell = 'E16'; Value = '  +0.26%  '; ForceText = $false }
    @{ Cell = 'D17'; Value = '0.0000171'; ForceText = $true }
    @{ Cell = 'E17'; Value = '  -1.63%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '2.527.64'; ForceText = $false }
    @{ Cell = 'E18'; Value = '  +3.78%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '11.02'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  -5.50%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '7.43'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  -6.13%  '; ForceText = $false }
    @{ Cell = 'E21'; Value = '  -3.50%  '; ForceText = $false }
    @{ Cell = 'D22'; Value = '4.02'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  -0.71%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '1.00'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  +0.31%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '4.25'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  -4.65%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '68.53'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  -3.31%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '1.80'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  -2.05%  '; ForceText = $false }
    @{ Cell = 'E27'; Value = '  -2.21%  '; ForceText = $false }
    @{ Cell = 'E28'; Value = '  +0.13%  '; ForceText = $false }
    @{ Cell = 'D29'; Value = '2.614.46'; ForceText = $false }
    @{ Cell = 'E29'; Value = '  +0.43%  '; ForceText = $false }
    @{ Cell = 'D30'; Value = '0.0₃0904'; ForceText = $false }
    @{ Cell = 'E30'; Value = '  -2.84%  '; ForceText = $false }
    @{ Cell = 'D31'; Value = '512.11'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  -0.93%  '; ForceText = $false }
    @{ Cell = 'D32'; Value = '7.76'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  -4.03%  '; ForceText = $false }
    @{ Cell = 'D33'; Value = '1.24'; ForceText = $true }
    @{ Cell = 'E33'; Value = '  -3.07%  '; ForceText = $false }
    @{ Cell = 'D34'; Value = '1.77'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  -3.26%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '0.999'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  -0.10%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '159.88'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  +1.03%  '; ForceText = $false }
    @{ Cell = 'E37'; Value = '  -6.88%  '; ForceText = $false }
    @{ Cell = 'D38'; Value = '18.71'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  +0.84%  '; ForceText = $false }
    @{ Cell = 'E39'; Value = '  -3.74%  '; ForceText = $false }
    @{ Cell = 'D40'; Value = '1.34'; ForceText = $true }
    @{ Cell = 'E40'; Value = '  -5.49%  '; ForceText = $false }
    @{ Cell = 'E41'; Value = '  -0.18%  '; ForceText = $false }
    @{ Cell = 'E42'; Value = '  -3.29%  '; ForceText = $false }
    @{ Cell = 'B43'; Value = 'RenderToken'; ForceText = $false }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'; ForceText = $false }
    @{ Cell = 'D43'; Value = '4.84'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  -2.45%  '; ForceText = $false }
    @{ Cell = 'B44'; Value = 'PolygonEcosystemToken'; ForceText = $false }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'; ForceText = $false }
    @{ Cell = 'D44'; Value = '0.328'; ForceText = $true }
    @{ Cell = 'E44'; Value = '  -1.56%  '; ForceText = $false }
    @{ Cell = 'E45'; Value = '  -2.53%  '; ForceText = $false }
    @{ Cell = 'E46'; Value = '  -1.18%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '143.05'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  +0.11%  '; ForceText = $false }
    @{ Cell = 'D48'; Value = '0.516'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  -4.03%  '; ForceText = $false }
    @{ Cell = 'D49'; Value = '3.46'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  -3.85%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '0.0₆0252'; ForceText = $false }
    @{ Cell = 'E50'; Value = '  -6.42%  '; ForceText = $false }
    @{ Cell = 'B51'; Value = 'Optimism'; ForceText = $false }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'; ForceText = $false }
    @{ Cell = 'D51'; Value = '1.57'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  -4.60%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Leading apostrophe forces Excel to store the literal text instead of
        # parsing it as a number (mirrors typing, e.g., `'25.50` into a cell).
        $range.Value = "'" + $u.Value
        # Re-normalize the style so the cell does not keep a distinct
        # quote-prefix style from the cells around it (it stays TEXT either way).
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
